$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1790.8636
$ws.Range("I40").Value = 1690.8182
$ws.Range("J40").Value = 1890.909
$ws.Range("K40").Value = 1690.8182
$ws.Range("L40").Value = 1890.909
$ws.Range("M40").Value = -1515.8182
$ws.Range("N40").Value = -2240.909

$ws.Range("H51").Value = 2183.2222
$ws.Range("I51").Value = 1188.8889
$ws.Range("J51").Value = 3177.5557
$ws.Range("K51").Value = 1188.8889
$ws.Range("L51").Value = 3177.5557
$ws.Range("M51").Value = -704.8888999999999
$ws.Range("N51").Value = -4145.5557

$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("N65").ClearContents()

$ws.Range("H106").Value = 252759.66
$ws.Range("I106").Value = 335546.2
$ws.Range("K106").Value = 335546.2
$ws.Range("M106").Value = -334915.2

$ws.Range("H116").Value = 3242.8333
$ws.Range("J116").Value = 3242.8333
$ws.Range("L116").Value = 3242.8333
$ws.Range("N116").Value = -10126.8333

$ws.Range("H137").Value = 6061525.5
$ws.Range("I137").Value = 896.3684
$ws.Range("J137").Value = 14286665
$ws.Range("K137").Value = 2689.1052
$ws.Range("L137").Value = 42859995
$ws.Range("M137").Value = -139.1052
$ws.Range("N137").Value = -42865095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2661.2
$ws.Range("I45").Value = 3042.4
$ws.Range("K45").Value = 3042.4
$ws.Range("M45").Value = -2665.4

$ws.Range("H64").Value = 22045.5
$ws.Range("J64").Value = 40091
$ws.Range("L64").Value = 40091
$ws.Range("N64").Value = -40587

$ws.Range("H67").Value = 22045.5
$ws.Range("J67").Value = 40091
$ws.Range("L67").Value = 40091
$ws.Range("N67").Value = -41807

$ws.Range("H74").Value = 15154292
$ws.Range("I74").Value = 23811068
$ws.Range("J74").Value = 4935.6665
$ws.Range("K74").Value = 23811068
$ws.Range("L74").Value = 4935.6665
$ws.Range("M74").Value = -23810194
$ws.Range("N74").Value = -6683.6665

$ws.Range("H77").Value = 15154292
$ws.Range("I77").Value = 23811068
$ws.Range("J77").Value = 4935.6665
$ws.Range("K77").Value = 119055340
$ws.Range("L77").Value = 24678.3325
$ws.Range("M77").Value = -119050972
$ws.Range("N77").Value = -33414.3325

$ws.Range("H132").Value = 6099594.5
$ws.Range("I132").Value = 8622476
$ws.Range("J132").Value = 2629.4167
$ws.Range("K132").Value = 25867428
$ws.Range("L132").Value = 7888.250100000001
$ws.Range("M132").Value = -25864898
$ws.Range("N132").Value = -12948.2501

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4648.5186
$ws.Range("I105").Value = 3955.4546
$ws.Range("J105").Value = 5125
$ws.Range("K105").Value = 3955.4546
$ws.Range("L105").Value = 5125
$ws.Range("M105").Value = -2208.4546
$ws.Range("N105").Value = -8619

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7096314.5
$ws.Range("I31").Value = 4746.647
$ws.Range("K31").Value = 4746.647
$ws.Range("M31").Value = -4451.647

$ws.Range("H34").Value = 7096314.5
$ws.Range("I34").Value = 4746.647
$ws.Range("K34").Value = 4746.647
$ws.Range("M34").Value = -4544.647

$ws.Range("H105").Value = 1428.8
$ws.Range("I105").Value = 1388
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 1388
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = 359
$ws.Range("N105").Value = -5494

$ws.Range("H140").Value = 44339.89
$ws.Range("J140").Value = 44339.89
$ws.Range("L140").Value = 44339.89
$ws.Range("N140").Value = -54699.89

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 1335.3334
$ws.Range("J44").Value = 1765.5
$ws.Range("L44").Value = 5296.5
$ws.Range("N44").Value = -6092.5

$ws.Range("H46").Value = 836.36365
$ws.Range("I46").Value = 425
$ws.Range("J46").Value = 1933.3334
$ws.Range("K46").Value = 1275
$ws.Range("L46").Value = 5800.0002
$ws.Range("M46").Value = -1184
$ws.Range("N46").Value = -5982.0002

$ws.Range("H50").Value = 205.07692
$ws.Range("I50").Value = 151.45454
$ws.Range("J50").Value = 500
$ws.Range("K50").Value = 454.36362
$ws.Range("L50").Value = 1500
$ws.Range("M50").Value = 26.63637999999997
$ws.Range("N50").Value = -2462

$ws.Range("H53").Value = 205.07692
$ws.Range("I53").Value = 151.45454
$ws.Range("J53").Value = 500
$ws.Range("K53").Value = 454.36362
$ws.Range("L53").Value = 1500
$ws.Range("M53").Value = 26.63637999999997
$ws.Range("N53").Value = -2462

$ws.Range("H98").Value = 1172.0834
$ws.Range("I98").Value = 650
$ws.Range("J98").Value = 1276.5
$ws.Range("K98").Value = 1950
$ws.Range("L98").Value = 3829.5
$ws.Range("M98").Value = -452
$ws.Range("N98").Value = -6825.5

$ws.Range("H131").Value = 859.28
$ws.Range("I131").Value = 790
$ws.Range("J131").Value = 859.9798
$ws.Range("K131").Value = 2370
$ws.Range("L131").Value = 2579.9394
$ws.Range("M131").Value = 2670
$ws.Range("N131").Value = -12659.9394

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2002.3462
$ws.Range("I102").Value = 1820.0435
$ws.Range("K102").Value = 1820.0435
$ws.Range("M102").Value = -198.0435

$ws.Range("H113").Value = 112510.336
$ws.Range("J113").Value = 1500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -5840

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 5359.364
$ws.Range("I122").Value = 6450.9165
$ws.Range("J122").Value = 4049.5
$ws.Range("K122").Value = 19352.7495
$ws.Range("L122").Value = 12148.5
$ws.Range("M122").Value = -16902.7495
$ws.Range("N122").Value = -17048.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1670.3636
$ws.Range("I100").Value = 2134
$ws.Range("J100").Value = 1496.5
$ws.Range("K100").Value = 4268
$ws.Range("L100").Value = 2993
$ws.Range("M100").Value = -3727
$ws.Range("N100").Value = -4075
